$d = $word.ActiveDocument

# 1. Merge "Graph Sandbox" + " Development Journal" into a single run.
$d.Content.Find.Execute("Graph Sandbox Development Journal", $false, $false, $false, $false, $false, $true, 1, $false, "Graph Sandbox Development Journal", 2) | Out-Null

# 2. Merge the curly-quoted epigraph runs into a single run.
$quote = [char]8220 + "To model interconnectedness" + [char]8230 + [char]8221
$d.Content.Find.Execute($quote, $false, $false, $false, $false, $false, $true, 1, $false, $quote, 2) | Out-Null

# 3. Merge "November 25" + ", 2020" into a single run.
$d.Content.Find.Execute("November 25, 2020", $false, $false, $false, $false, $false, $true, 1, $false, "November 25, 2020", 2) | Out-Null

# 4. Insert a new paragraph before the "BinaryTree, AVLTree..." paragraph and
#    move the _GoBack bookmark from the end of the "Somehow the algorithms..."
#    paragraph to the end of this new paragraph.
$target = $d.Paragraphs.Item(12)
$target.Range.InsertParagraphBefore() | Out-Null
$newPara = $d.Paragraphs.Item(12)
$newPara.Range.Text = "We want support for directed and undirected graphs/edges too. I" + [char]8217 + "m starting to think it might be better to keep the tree and graph projects separate, or at least keep their code fairly separate in this project. Trees and graphs don" + [char]8217 + "t have as much in common as I had thought. "

$newPara = $d.Paragraphs.Item(12)
$bookmarkRange = $d.Range($newPara.Range.End - 1, $newPara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
